$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.883.87'
$ws.Range("E2").Value = '  -0.02%  '

# Row 3: Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.639.36'
$ws.Range("E3").Value = '  -0.12%  '

# Row 4: TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  -0.19%  '

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.20'
$ws.Range("E5").Value = '  -0.31%  '

# Row 6: XRP
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5066'
$ws.Range("E6").Value = '  +0.55%  '

# Row 7: USDC
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.004'
$ws.Range("E7").Value = '  -0.12%  '

# Row 8: Cardano
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2561'
$ws.Range("E8").Value = '  -0.47%  '

# Row 9: Dogecoin
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06362'
$ws.Range("E9").Value = '  -0.30%  '

# Row 10: Solana
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.47'
$ws.Range("E10").Value = '  -1.04%  '

# Row 11: TRON
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07767'
$ws.Range("E11").Value = '  -0.08%  '

# Row 12: Polkadot
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.277'
$ws.Range("E12").Value = '  +0.39%  '

# Row 13: WrappedEther
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.645.51'
$ws.Range("E13").Value = '  +0.24%  '

# Row 14: Polygon
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5440'
$ws.Range("E14").Value = '  -0.54%  '

# Row 15: ShibaInu
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0₅7811'
$ws.Range("E15").Value = '  -1.13%  '

# Row 16: Litecoin
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.29'
$ws.Range("E16").Value = '  +0.27%  '

# Row 17: WrappedBTC
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '25.956.72'
$ws.Range("E17").Value = '  +0.20%  '

# Row 18: Dai
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.004'
$ws.Range("E18").Value = '  -0.10%  '

# Row 19: BitcoinCash
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '197.34'
$ws.Range("E19").Value = '  -2.68%  '

# Row 20: Uniswap
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.450'
$ws.Range("E20").Value = '  +1.33%  '

# Row 21: Avalanche
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.934'
$ws.Range("E21").Value = '  +0.21%  '

# Row 22: Chainlink
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.037'
$ws.Range("E22").Value = '  +0.92%  '

# Row 23: BinanceUSD
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.004'
$ws.Range("E23").Value = '  -0.19%  '

# Row 24: Toncoin
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.882'
$ws.Range("E24").Value = '  -1.50%  '

# Row 25: Monero
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '141.17'
$ws.Range("E25").Value = '  +0.15%  '

# Row 26: Stellar
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1164'
$ws.Range("E26").Value = '  +2.72%  '

# Row 27: Cosmos
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.868'
$ws.Range("E27").Value = '  +1.26%  '

# Row 28: EthereumClassic
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.70'
$ws.Range("E28").Value = '  +0.09%  '

# Row 29: PancakeSwap
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.238'
$ws.Range("E29").Value = '  -0.55%  '

# Row 30: Hedera
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.04995'
$ws.Range("E30").Value = '  +0.26%  '

# Row 31: InternetComputer(DFINITY)
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.251'
$ws.Range("E31").Value = '  -0.55%  '

# Row 32: Filecoin
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.187'
$ws.Range("E32").Value = '  -0.36%  '

# Row 33: LidoDAOToken
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.540'
$ws.Range("E33").Value = '  -0.57%  '

# Row 34: HuobiToken
$ws.Range("E34").Value = '  -0.47%  '

# Row 35: ARBITRUM
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.8953'
$ws.Range("E35").Value = '  +0.14%  '

# Row 36: MXToken
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.584'
$ws.Range("E36").Value = '  -1.91%  '

# Row 37: Maker
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.126.86'
$ws.Range("E37").Value = '  -1.86%  '

# Row 38: ImmutableX
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5456'
$ws.Range("E38").Value = '  -3.09%  '

# Row 39: VeChain
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01557'
$ws.Range("E39").Value = '  -0.73%  '

# Row 40: PaxDollar
$ws.Range("E40").Value = '  -0.16%  '

# Row 41: BabyDogeCoin
$ws.Range("B41").Value = 'BabyDogeCoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0₈131'
$ws.Range("E41").Value = '  +12.19%  '

# Row 42: mCoin
$ws.Range("B42").Value = 'mCoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.544'
$ws.Range("E42").Value = '  -1.20%  '

# Row 43: FraxShare
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.592'
$ws.Range("E43").Value = '  -1.63%  '

# Row 44: TrustWalletToken
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8171'
$ws.Range("E44").Value = '  +1.08%  '

# Row 45: Quant
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '99.79'
$ws.Range("E45").Value = '  -0.28%  '

# Row 46: RocketPoolETH
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.783.87'
$ws.Range("E46").Value = '  +0.42%  '

# Row 47: Mantle
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4538'
$ws.Range("E47").Value = '  -0.18%  '

# Row 48: Frax
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.004'
$ws.Range("E48").Value = '  -0.29%  '

# Row 49: Aave
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '54.73'
$ws.Range("E49").Value = '  -0.36%  '

# Row 50: Cronos
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05070'
$ws.Range("E50").Value = '  +0.23%  '

# Row 51: USDD
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.005'
$ws.Range("E51").Value = '  +0.22%  '

Write-Output "applied 106 cell updates"
